# Generate Report for Handback
# Re-run of the handback status report: two file GUIDs were regenerated,
# their content hashes changed, and the handoff/handback timestamps moved
# forward ~1 minute. This updates every cell value and hyperlink display
# text that mirrors the old identifiers/timestamps to the new ones, across
# the Overview / zh-cn / de-de sheets.

$map = @{
    "11f943d4-e580-437c-8da3-f377f2e05c15.md" = "56426c72-e1a0-49f7-bb3a-df81ef7ff23f.md";
    "bb898a2c-ec1c-4f15-987d-4e981b05dd91.md" = "ffffcf2e933f-cbf8-4297-a5f3-fbd5915e14d0.md";
    "11f943d4-e580-437c-8da3-f377f2e05c15.3f056b18b7de5a52b2f2b3434c918479e351e156.zh-cn.xlf" = "56426c72-e1a0-49f7-bb3a-df81ef7ff23f.8d2f12b505b9981efa8ba5b178d3047ae1a44b04.zh-cn.xlf";
    "bb898a2c-ec1c-4f15-987d-4e981b05dd91.7f7a0c4fed9b66d3029f48ed7d568ae715cbcbe9.zh-cn.xlf" = "56426c72-e1a0-49f7-bb3a-df81ef7ff23f.8d2f12b505b9981efa8ba5b178d3047ae1a44b04.zh-cn.xlf";
    "11f943d4-e580-437c-8da3-f377f2e05c15.3f056b18b7de5a52b2f2b3434c918479e351e156.de-de.xlf" = "56426c72-e1a0-49f7-bb3a-df81ef7ff23f.8d2f12b505b9981efa8ba5b178d3047ae1a44b04.de-de.xlf";
    "bb898a2c-ec1c-4f15-987d-4e981b05dd91.7f7a0c4fed9b66d3029f48ed7d568ae715cbcbe9.de-de.xlf" = "56426c72-e1a0-49f7-bb3a-df81ef7ff23f.8d2f12b505b9981efa8ba5b178d3047ae1a44b04.de-de.xlf";
    "2016-03-18 22:49:25" = "2016-03-18 22:50:32";
    "2016-03-18 22:49:44" = "2016-03-18 22:50:52";
    "2016-03-18 22:49:28" = "2016-03-18 22:50:35";
    "2016-03-18 22:49:49" = "2016-03-18 22:50:57";
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Update plain cell values (covers columns without hyperlinks, e.g.
    # the datetime columns, as well as the text underneath hyperlinked
    # cells so the shared-string table tracks the display text).
    # NOTE: read via .Text (not .Value - the getter does not resolve
    # through this host), write via .Value.
    $used = $ws.UsedRange
    foreach ($row in 1..$used.Rows.Count) {
        foreach ($col in 1..$used.Columns.Count) {
            $cell = $ws.Cells.Item($row, $col)
            $old = $cell.Text
            if ($old -ne $null -and $map.ContainsKey($old)) {
                $cell.Value = $map[$old]
            }
        }
    }

    # Update the hyperlink display text to match (in place, via the
    # enumerator - Item() on this collection does not resolve to an
    # individual link, so iterate with foreach instead).
    foreach ($hl in $ws.Hyperlinks) {
        $oldDisplay = $hl.TextToDisplay
        if ($oldDisplay -ne $null -and $map.ContainsKey($oldDisplay)) {
            $hl.TextToDisplay = $map[$oldDisplay]
        }
    }
}
